$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 14 content/formatting update ------------------------------------
# G14: add a new note describing the entity-count behaviour for projectiles.
$ws.Cells.Item(14, 7).Value = "When a projectile is fired - increment entity, when an object falls out of octree, decrement"

# H14: move status from "TODO" (red, no border) to "UNDERWAY" (theme accent2
# fill, no border). Clone the fill used by the other "UNDERWAY" cells (H3)
# and then strip the border that comes along with that formatting, leaving
# just the fill - matching the rest of the non-header status cells.
$h3 = $ws.Cells.Item(3, 8)
$h14 = $ws.Cells.Item(14, 8)
$h3.Copy()
$h14.PasteSpecial(-4122)   # xlPasteFormats
$h14.Borders.LineStyle = -4142   # xlLineStyleNone
$h14.Value = "UNDERWAY"

# --- View selection ---------------------------------------------------------
$ws.Range("I28").Select() | Out-Null

$excel.CutCopyMode = 0
